$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit rotates rows 10-12 (columns A,B,D,E,F,G,H,Q,R,AC):
#   old row12 -> row10, old row10 -> row11, old row11 -> row12

$cols = @(1,2,4,5,6,7,8,17,18,29)

$row10 = @{}
$row11 = @{}
$row12 = @{}
foreach ($c in $cols) {
    $row10[$c] = $ws.Cells.Item(10, $c).Value()
    $row11[$c] = $ws.Cells.Item(11, $c).Value()
    $row12[$c] = $ws.Cells.Item(12, $c).Value()
}

foreach ($c in $cols) {
    $ws.Cells.Item(10, $c).Value = $row12[$c]
    $ws.Cells.Item(11, $c).Value = $row10[$c]
    $ws.Cells.Item(12, $c).Value = $row11[$c]
}
